$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "vat" column: header in M1 plus values in M2:M4, matching the
# existing "Times New Roman 14, wrap text, vertical centered" style (s="1")
# used throughout the rest of the table.
$vatRange = $ws.Range("M1:M4")
$vatRange.Font.Name = "Times New Roman"
$vatRange.Font.Size = 14
$vatRange.WrapText = $true
$vatRange.VerticalAlignment = -4108

$ws.Range("M1").Value = "vat"
$ws.Range("M2").Value = 5
$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 2

# Update the view: scroll so column B is left-most visible, select N2 (as in
# the target workbook's sheetView).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("N2").Select()
